$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format before assigning values, so that
# numeric-looking price strings (e.g. "1.00", "43.60") are preserved exactly
# as text instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.922.44'
$ws.Range("E2").Value = '  +1.63%  '
$ws.Range("D3").Value = '3.803.08'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.36%  '
$ws.Range("D5").Value = '445.00'
$ws.Range("E5").Value = '  +7.16%  '
$ws.Range("D6").Value = '146.42'
$ws.Range("E6").Value = '  +15.95%  '
$ws.Range("D7").Value = '0.622'
$ws.Range("E7").Value = '  +5.13%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").Value = '0.737'
$ws.Range("E9").Value = '  +4.14%  '
$ws.Range("E10").Value = '  -1.47%  '
$ws.Range("D11").Value = '0.0000320'
$ws.Range("E11").Value = '  -5.57%  '
$ws.Range("D12").Value = '43.60'
$ws.Range("E12").Value = '  +11.43%  '
$ws.Range("D13").Value = '10.30'
$ws.Range("E13").Value = '  +4.57%  '
$ws.Range("D14").Value = '4.409.26'
$ws.Range("E14").Value = '  +1.07%  '
$ws.Range("D15").Value = '14.89'
$ws.Range("E15").Value = '  -8.26%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.811.23'
$ws.Range("E16").Value = '  +1.49%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = '0.137'
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").Value = '19.95'
$ws.Range("E18").Value = '  +4.33%  '
$ws.Range("D19").Value = '1.14'
$ws.Range("E19").Value = '  +8.20%  '
$ws.Range("D20").Value = '67.035.64'
$ws.Range("E20").Value = '  +1.47%  '
$ws.Range("D21").Value = '422.85'
$ws.Range("E21").Value = '  +5.93%  '
$ws.Range("D22").Value = '14.57'
$ws.Range("E22").Value = '  +4.53%  '
$ws.Range("D23").Value = '3.25'
$ws.Range("E23").Value = '  +10.77%  '
$ws.Range("D24").Value = '86.36'
$ws.Range("E24").Value = '  +4.78%  '
$ws.Range("D25").Value = '37.19'
$ws.Range("E25").Value = '  +2.80%  '
$ws.Range("E26").Value = '  +8.79%  '
$ws.Range("E27").Value = '  -0.21%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '9.48'
$ws.Range("E28").Value = '  +22.24%  '
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").Value = '9.71'
$ws.Range("E29").Value = '  +6.36%  '
$ws.Range("D30").Value = '737.51'
$ws.Range("E30").Value = '  +6.40%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = '0.133'
$ws.Range("E31").Value = '  +13.17%  '
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").Value = '13.62'
$ws.Range("E32").Value = '  +13.02%  '
$ws.Range("D33").Value = '2.73'
$ws.Range("E33").Value = '  -0.69%  '
$ws.Range("D34").Value = '43.10'
$ws.Range("E34").Value = '  +17.17%  '
$ws.Range("D35").Value = '0.157'
$ws.Range("E35").Value = '  +7.42%  '
$ws.Range("D36").Value = '56.99'
$ws.Range("E36").Value = '  +4.94%  '
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("D38").Value = '5.48'
$ws.Range("E38").Value = '  +21.03%  '
$ws.Range("D39").Value = '0.0474'
$ws.Range("E39").Value = '  +6.77%  '
$ws.Range("D40").Value = '2.89'
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").Value = '0.341'
$ws.Range("E41").Value = '  +18.92%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.29%  '
$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").Value = '0.0₃0672'
$ws.Range("E43").Value = '  -9.50%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").Value = '0.140'
$ws.Range("E44").Value = '  +6.02%  '
$ws.Range("D45").Value = '2.51'
$ws.Range("E45").Value = '  +17.55%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '3.26'
$ws.Range("E46").Value = '  +5.94%  '
$ws.Range("B47").Value = 'LidoDAOToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D47").Value = '3.41'
$ws.Range("E47").Value = '  +4.69%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = '146.57'
$ws.Range("E48").Value = '  +2.81%  '
$ws.Range("D49").Value = '2.10'
$ws.Range("E49").Value = '  +5.31%  '
$ws.Range("B50").Value = 'WEMIXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").Value = '2.67'
$ws.Range("E50").Value = '  +7.45%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").Value = '2.86'
$ws.Range("E51").Value = '  +6.54%  '

# Restore the original (default) cell style so no new style/number format
# is left behind on the cells.
$ws.Range("D2:D51").Style = "Normal"
